$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Diplomado")
$ws2 = $wb.Worksheets.Item("Bachillerato")

# --- Sheet "Bachillerato" (sheet2) ---
$ws2.Range("C17").Value = "Total de créditos"
$ws2.Range("C1").Copy()
$ws2.Range("C17").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

$ws2.Range("D17").Formula = "=SUM(D2:D16)"
$ws2.Range("D2").Copy()
$ws2.Range("D17").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false

[void]$ws2.Range("E14").Select()

# --- Sheet "Diplomado" (sheet1) ---
# J7 gets a value of 10 (row for "Principios de Administracion")
$ws1.Range("J7").Value = 10

# New total row (row 25): label in D25 styled like the bold "Total" header
# (white bold text on black fill, left/right borders only)
$ws1.Range("D25").Value = "Total de créditos"
$ws1.Range("A1").Copy()
$ws1.Range("D25").PasteSpecial(-4122)
$ws1.Range("D25").Borders.Item(8).LineStyle = -4142
$ws1.Range("D25").Borders.Item(9).LineStyle = -4142
$ws1.Application.CutCopyMode = $false

# E25 sum formula, styled like the other bordered/centered data cells
$ws1.Range("E25").Formula = "=SUM(E2:E24)"
$ws2.Range("A2").Copy()
$ws1.Range("E25").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = $false

[void]$ws1.Activate()
[void]$ws1.Range("J8").Select()
